$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Final day (Day 11 / column O) effort numbers: both tasks logged 5 hours
# of work on the last day of the sprint, so "Completed Effort" (row 8,
# the SUM formula) rises to 10 and "Remaining Effort" (row 9) drops to 0.
$ws.Range("O6").Value = 5
$ws.Range("O7").Value = 5
